# Refactor synthetic array: insert a new "statut_name" column (C) that
# spells out, in French, what each "statut_label" code (noir/rouge/
# orange/vert) means. Everything that used to live in columns C..L
# (NCTId .. intervention_type) shifts one column to the right (D..M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the old C:L block one column to the right by inserting a fresh,
# blank column at C - this is exactly what Excel's own
# "Insert Sheet Column" does and it carries all existing data/styles
# with it automatically.
$ws.Columns("C:C").Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 3).Value = "statut_name"

# Human-readable label, keyed off the existing statut_label (col B).
$labelMap = @{
    "noir"   = "pas de résultat ni de publication";
    "rouge"  = "résultat et / ou publication posté";
    "orange" = "résultat et / ou publication posté dans les 36 mois";
    "vert"   = "résultat et / ou publication posté dans les 12 mois";
}

$lastRow = $ws.Cells.Item(1, 1).End(-4121).Row()  # xlDown
for ($r = 2; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($r, 3).Value = $labelMap[$label]
}
